$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values in column D look like plain decimal numbers (e.g. "235.64").
# Excel's automatic type detection would otherwise silently convert these into
# numeric cells (losing the exact textual representation, e.g. "58.00" -> 58).
# The source workbook stores every Price cell as text, so force a Text number
# format on those specific cells before writing the value, to keep them as text.
$textPriceRows = @(5,10,14,15,16,17,20,21,23,25,26,27,28,29,30,31,34,35,36,38,39,41,44,45,48,50)
foreach ($r in $textPriceRows) {
  $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.196.04"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.053.80"
$ws.Range("E3").Value = "  +3.92%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "235.64"
$ws.Range("E5").Value = "  -2.30%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.94%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - Solana
$ws.Range("E8").Value = "  +4.76%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.65%  "

# Row 10 - OKB
$ws.Range("D10").Value = "58.00"
$ws.Range("E10").Value = "  -1.43%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.17%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +3.05%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.356.22"
$ws.Range("E13").Value = "  +3.92%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "14.54"
$ws.Range("E14").Value = "  +4.60%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "20.78"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.773"
$ws.Range("E16").Value = "  +3.40%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "5.22"
$ws.Range("E17").Value = "  +4.48%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.067.66"
$ws.Range("E18").Value = "  +4.72%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "37.321.19"
$ws.Range("E19").Value = "  +1.12%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "5.90"
$ws.Range("E20").Value = "  +19.81%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "68.29"
$ws.Range("E21").Value = "  +0.50%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "223.04"
$ws.Range("E23").Value = "  -1.67%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.15%  "

# Row 25 - now PancakeSwap (was Toncoin)
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  +3.86%  "

# Row 26 - now Toncoin (was PancakeSwap)
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  +3.41%  "

# Row 27 - Monero
$ws.Range("D27").Value = "163.06"
$ws.Range("E27").Value = "  +1.09%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "8.86"
$ws.Range("E28").Value = "  +3.58%  "

# Row 29 - Kaspa
$ws.Range("D29").Value = "0.131"
$ws.Range("E29").Value = "  +4.88%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "19.20"
$ws.Range("E30").Value = "  +1.26%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "1.37"
$ws.Range("E31").Value = "  +7.29%  "

# Row 32 - Stellar
$ws.Range("E32").Value = "  +0.91%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +1.87%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.0617"
$ws.Range("E34").Value = "  +1.31%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "2.51"
$ws.Range("E35").Value = "  +8.71%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("D36").Value = "4.34"
$ws.Range("E36").Value = "  +2.89%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  +0.04%  "

# Row 38 - THORChain
$ws.Range("D38").Value = "5.93"
$ws.Range("E38").Value = "  +16.35%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "3.32"
$ws.Range("E39").Value = "  +1.95%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  -1.41%  "

# Row 41 - HuobiToken
$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  -1.98%  "

# Row 42 - Cronos
$ws.Range("E42").Value = "  +7.44%  "

# Row 43 - now Maker (was FTXToken)
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.482.42"
$ws.Range("E43").Value = "  +5.08%  "

# Row 44 - now FTXToken (was Maker)
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "4.32"
$ws.Range("E44").Value = "  +10.43%  "

# Row 45 - Aave
$ws.Range("D45").Value = "94.36"
$ws.Range("E45").Value = "  +8.84%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +2.55%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +1.27%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "15.98"
$ws.Range("E48").Value = "  +5.42%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +3.42%  "

# Row 50 - FraxShare
$ws.Range("D50").Value = "7.13"
$ws.Range("E50").Value = "  +8.41%  "

# Row 51 - MXToken
$ws.Range("E51").Value = "  +2.28%  "
